# Add exclusion plots for Co1
#
# For the six non-CO1 result sheets (Sheet2, Sheet3, Sheet4, Sheet5, Sheet7,
# Sheet8) a new "CustomNBC" row is inserted immediately above the
# "Kraken2_0.0" row (row 4). Its values are copies of what the existing
# "NBC" row (row 10) already holds. The old "NBC" row (now shifted down to
# row 11 after the insert) is then removed, so the net effect is that
# Kraken2_0.0 .. Mothur each move down by one row and a brand-new
# "CustomNBC" row takes their old starting slot.

$wb = $excel.ActiveWorkbook

# sheetName -> (B, C, D, E, F, G, H) values copied from that sheet's
# pre-existing "NBC" row (row 10).
$sheetData = @{
    "Sheet2" = @("100 Australian species", "12S", 0.6885245901639344, 0.525, 0.5957446808510638, 0.6481481481481481, 0.4242424242424243)
    "Sheet3" = @("100 Australian species", "16S", 0.8387096774193549, 0.5842696629213483, 0.6887417218543046, 0.7715133531157269, 0.5252525252525253)
    "Sheet4" = @("Lutjanidae", "12S", 0.8, 0.631578947368421, 0.7058823529411765, 0.759493670886076, 0.5833333333333334)
    "Sheet5" = @("Lutjanidae", "16S", 0.9, 0.36, 0.5142857142857143, 0.6923076923076923, 0.3703703703703703)
    "Sheet7" = @("Rottnest", "12S", 0.7571428571428571, 0.6309523809523809, 0.6883116883116883, 0.7280219780219781, 0.5294117647058824)
    "Sheet8" = @("Rottnest", "16S", 0.7671232876712328, 0.6153846153846154, 0.6829268292682927, 0.7310704960835509, 0.5357142857142857)
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $sheetData[$sheetName]

    # Insert a fresh row at 4 -- this pushes the existing Kraken2_0.0..NBC
    # block (rows 4-10) down to rows 5-11 without touching their stored
    # values (full precision preserved).
    $ws.Rows.Item(4).Insert()

    # Populate the new row 4 with the CustomNBC entry.
    $ws.Cells.Item(4, 1).Value = "CustomNBC"
    $ws.Cells.Item(4, 2).Value = $vals[0]
    $ws.Cells.Item(4, 3).Value = $vals[1]
    $ws.Cells.Item(4, 4).Value = $vals[2]
    $ws.Cells.Item(4, 5).Value = $vals[3]
    $ws.Cells.Item(4, 6).Value = $vals[4]
    $ws.Cells.Item(4, 7).Value = $vals[5]
    $ws.Cells.Item(4, 8).Value = $vals[6]

    # The original "NBC" row is now at row 11 -- remove it so Mothur ends
    # up as the last row of the tool block (row 10), matching the target
    # layout.
    $ws.Rows.Item(11).Delete()
}

"done"
